# Allowed Hosts workbook maintenance:
#  - prune the host list down to the currently-valid entries
#  - rename the "Maryland Comcast" nickname to "Priscilla's"
#  - tidy selection / formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows that no longer correspond to an allowed host, working from
# the bottom up so earlier row numbers stay valid as rows shift up.
$rowsToRemove = @(25,24,23,21,20,18,17,16,15,14,10,8,7,6,5,4,3,2)
foreach ($r in $rowsToRemove) {
    $ws.Rows($r).Delete()
}

# The row that used to be labelled "Maryland Comcast" survives (it's now
# row 7) but the nickname has changed.
$ws.Range("A7").Value = "Priscilla's"

# Give the surviving data rows an explicit (Normal) font so the sheet's
# formatting is no longer purely inherited.
$ws.Range("A2:B7").Font.ThemeColor = 1

# Move the active cell back to column A and refresh the used range.
$ws.Range("A9").Select()
